# Remove the "Lan dau" (first-round) approval columns G:H entirely. This
# shifts the former "Lan cuoi" (final-round) columns (old I,J) left into
# G,H and the "Ghi chu" column (old K) left into I - matching column
# widths, dimensions and shared-string usage automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1:H1").EntireColumn.Delete()

# Rename the surviving "Lan cuoi" headers, dropping the now-redundant
# "Lan cuoi" qualifier since there is only one approval stage left.
$ws.Range("G1").Value = "Ngày Duyệt/Từ chối"
$ws.Range("H1").Value = "Trạng thái"

# Consolidate the redundant left-aligned/wrap body style onto the
# identical pre-existing style so these cells stop using the duplicate
# cell format.
$ws.Range("C2").HorizontalAlignment = -4131
$ws.Range("E2").HorizontalAlignment = -4131
$ws.Range("G2").HorizontalAlignment = -4131
$ws.Range("H2").HorizontalAlignment = -4131
$ws.Range("I2").HorizontalAlignment = -4131

# Update the saved view/selection to match the edited layout.
$ws.Range("H2").Select()
